# Logged Week 15 and simulated Week 16
# Update the "H" row target depth data on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# Offense sheet ("OFF") - row 2 is the "H" (Home) row
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 214
$wsOff.Range("C2").Value = 142
$wsOff.Range("D2").Value = 51
$wsOff.Range("E2").Value = 18
$wsOff.Range("F2").Value = 3

# Defense sheet ("DEF") - row 2 is the "H" (Home) row
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 252
$wsDef.Range("C2").Value = 176
$wsDef.Range("D2").Value = 49
$wsDef.Range("E2").Value = 24
$wsDef.Range("F2").Value = 7
